$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

# --- Row 55: 2020-07-24 -------------------------------------------------
# Temporarily format column A as text so the ISO-looking date string is
# stored as a literal shared string (matching the rest of column A) rather
# than being auto-converted into a date serial number, then restore the
# default style so no visible formatting change is introduced.
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "2020-07-24"
$ws.Range("A55").Style = "Normal"

$ws.Range("B55").Value = 378285
$ws.Range("C55").Value = 426869
$ws.Range("D55").Value = 90970
$ws.Range("E55").Value = 42645
$ws.Range("F55").Value = 27.97

# --- Row 56: 2020-07-25 -------------------------------------------------
$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = "2020-07-25"
$ws.Range("A56").Style = "Normal"

$ws.Range("B56").Value = 385036
$ws.Range("C56").Value = 433384
$ws.Range("D56").Value = 93433
$ws.Range("E56").Value = 43374
$ws.Range("F56").Value = 27.89
